$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Equity tickers (column B)
$ws.Range("B4").Value = "EP"
$ws.Range("B5").Value = "ISRG"
$ws.Range("B6").Value = "CMD"
$ws.Range("B7").Value = "OKE"

# Prices (column C)
$ws.Range("C4").Value = 13.72
$ws.Range("C5").Value = 269
$ws.Range("C6").Value = 47.55
$ws.Range("C7").Value = 26.76

# Percentages (column D) - update D6/D7
$ws.Range("D6").Value = 0.3
$ws.Range("D7").Value = 0.5

# Formulas (column E), use absolute reference to C2
$ws.Range("E4").Formula = "=`$C`$2*D4/C4"
$ws.Range("E5").Formula = "=`$C`$2*D5/C5"
$ws.Range("E6").Formula = "=`$C`$2*D6/C6"
$ws.Range("E7").Formula = "=`$C`$2*D7/C7"

# Number formats
$ws.Range("C4:C7").Style = "Currency"
$ws.Range("E4:E7").NumberFormat = "0.00"

# Selection
$ws.Range("E4").Select() | Out-Null
